$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.368.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.813.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.48%  "

$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3812"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3495"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.81"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.04%  "

$ws.Range("E10").Value = "  +3.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07724"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.48%  "

$ws.Range("E12").Value = "  +0.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.02%  "

$ws.Range("E14").Value = "  +5.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.249"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.807.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001118"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06717"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.81%  "

$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.549"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.388.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.467"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.671"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.479"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.006.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "136.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.24%  "

$ws.Range("E32").Value = "  +6.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.044"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08739"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.61%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.708"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.609"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6973"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2277"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02416"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06496"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.947"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.297"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6512"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.017"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.177"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07346"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.56%  "
